$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The narrower Range object returned directly by Table.Cell(...).Range behaves
# oddly with Find/replace in this runtime (it can resolve to the wrong part of
# the document). To stay safe we only ever use it to read .Start/.Text, then
# build a brand-new plain $d.Range(start,end) for the actual text surgery.

function Set-CellSubstring($table, $row, $col, $search, $replacement) {
    $cell = $table.Cell($row, $col)
    $cellRange = $cell.Range
    $cellStart = $cellRange.Start
    $cellText = $cellRange.Text

    $idx = $cellText.IndexOf($search)
    if ($idx -lt 0) {
        Write-Host "NOT FOUND: row=$row col=$col search=[$search]"
        return
    }

    $absStart = $cellStart + $idx
    $absEnd = $absStart + $search.Length
    $target = $d.Range($absStart, $absEnd)

    if ($target.Text -ne $search) {
        Write-Host "MISMATCH: row=$row col=$col expected=[$search] got=[$($target.Text)]"
        return
    }

    $target.Text = $replacement
    Write-Host "OK: row=$row col=$col [$search] -> [$replacement]"
}

# 1. Quantity (line item): 1 -> 2
Set-CellSubstring $t 5 4 "1" "2"

# 2. Amount (line item): 30000.0 -> 60000.0
Set-CellSubstring $t 5 11 "30000.0" "60000.0"

# 3. Total quantity: 01 -> 02
Set-CellSubstring $t 6 4 "01" "02"

# 4. Total amount: 30000.0 -> 60000.0
Set-CellSubstring $t 6 11 "30000.0" "60000.0"

# 5. Amount in words: Thirty Thousand only -> Sixty Thousand only
Set-CellSubstring $t 7 1 "Thirty Thousand only" "Sixty Thousand only"

# 6. Amounts: Total figure: 30000.0 -> 60000.0
Set-CellSubstring $t 7 6 "30000.0" "60000.0"

# 7. Amounts: Received figure: 30000 -> 10000
Set-CellSubstring $t 8 6 "30000" "10000"
